$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.536.70"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.414.45"
$ws.Range("E3").Value = "  +8.69%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.76"
$ws.Range("E5").Value = "  +12.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.22"
$ws.Range("E6").Value = "  -5.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +3.53%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.655"
$ws.Range("E9").Value = "  +9.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.06"
$ws.Range("E10").Value = "  -3.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.72"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.03"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.35"
$ws.Range("E14").Value = "  +16.69%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.775.80"
$ws.Range("E16").Value = "  +8.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.499.57"
$ws.Range("E17").Value = "  +12.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.545.18"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("E20").Value = "  +3.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.46"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "261.56"
$ws.Range("E23").Value = "  +13.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.63"
$ws.Range("E25").Value = "  +7.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +5.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.95"
$ws.Range("E28").Value = "  +9.81%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "179.22"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.35"
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0934"
$ws.Range("E33").Value = "  +6.71%  "
$ws.Range("E34").Value = "  +7.03%  "
$ws.Range("E35").Value = "  +4.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.93"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.94"
$ws.Range("E38").Value = "  -5.93%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  +22.87%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.106"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.62"
$ws.Range("E41").Value = "  +24.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.235"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.93"
$ws.Range("E43").Value = "  +21.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.25"
$ws.Range("E44").Value = "  -7.59%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.67"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.69"
$ws.Range("E47").Value = "  +5.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.53"
$ws.Range("E48").Value = "  +13.30%  "
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.595.02"
$ws.Range("E50").Value = "  +12.96%  "
$ws.Range("E51").Value = "  +3.68%  "
